$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value while forcing text storage (avoids Excel
# auto-converting numeric-looking strings like "1.00" or "212.50" to numbers),
# then restores the cell to the default "Normal" style so no stray
# number-format style gets attached to the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '27.946.09'
$ws.Range("E2").Value = '  +1.70%  '

Set-TextValue $ws.Range("D3") '1.640.96'
$ws.Range("E3").Value = '  +0.14%  '

Set-TextValue $ws.Range("D4") '1.00'
$ws.Range("E4").Value = '  -0.18%  '

Set-TextValue $ws.Range("D5") '212.50'
$ws.Range("E5").Value = '  -0.05%  '

Set-TextValue $ws.Range("D6") '0.526'
$ws.Range("E6").Value = '  -1.62%  '

Set-TextValue $ws.Range("D7") '1.00'
$ws.Range("E7").Value = '  -0.23%  '

Set-TextValue $ws.Range("D8") '23.40'
$ws.Range("E8").Value = '  +1.22%  '

Set-TextValue $ws.Range("D9") '0.263'
$ws.Range("E9").Value = '  +2.80%  '

$ws.Range("E10").Value = '  +0.39%  '

Set-TextValue $ws.Range("D11") '0.0892'
$ws.Range("E11").Value = '  +0.63%  '

Set-TextValue $ws.Range("D12") '1.873.48'
$ws.Range("E12").Value = '  +0.15%  '

Set-TextValue $ws.Range("D13") '1.636.62'
$ws.Range("E13").Value = '  -0.27%  '

$ws.Range("E14").Value = '  +1.00%  '

$ws.Range("E15").Value = '  -3.79%  '

Set-TextValue $ws.Range("D16") '64.66'
$ws.Range("E16").Value = '  +0.82%  '

Set-TextValue $ws.Range("D17") '27.926.69'
$ws.Range("E17").Value = '  +1.74%  '

Set-TextValue $ws.Range("D18") '233.34'
$ws.Range("E18").Value = '  +1.75%  '

$ws.Range("E19").Value = '  +0.13%  '

Set-TextValue $ws.Range("D20") '7.64'
$ws.Range("E20").Value = '  +1.52%  '

$ws.Range("E21").Value = '  -0.26%  '

$ws.Range("E22").Value = '  +0.13%  '

Set-TextValue $ws.Range("D23") '10.00'
$ws.Range("E23").Value = '  +3.21%  '

$ws.Range("E24").Value = '  +4.33%  '

Set-TextValue $ws.Range("D25") '150.70'
$ws.Range("E25").Value = '  +1.01%  '

Set-TextValue $ws.Range("D26") '6.94'
$ws.Range("E26").Value = '  -0.61%  '

$ws.Range("E27").Value = '  -0.81%  '

Set-TextValue $ws.Range("D28") '15.70'
$ws.Range("E28").Value = '  +1.02%  '

$ws.Range("E29").Value = '  -0.19%  '

$ws.Range("E30").Value = '  +0.25%  '

Set-TextValue $ws.Range("D31") '0.0483'
$ws.Range("E31").Value = '  -0.61%  '

$ws.Range("E32").Value = '  +0.62%  '

Set-TextValue $ws.Range("D33") '1.474.22'
$ws.Range("E33").Value = '  +4.03%  '

$ws.Range("E34").Value = '  -1.81%  '

$ws.Range("E35").Value = '  -2.55%  '

$ws.Range("E36").Value = '  -0.52%  '

$ws.Range("E37").Value = '  -0.37%  '

$ws.Range("E38").Value = '  +0.05%  '

$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D39") '0.0168'
$ws.Range("E39").Value = '  +0.67%  '

$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range("D40") '0.921'
$ws.Range("E40").Value = '  +13.30%  '

Set-TextValue $ws.Range("D41") '69.20'
$ws.Range("E41").Value = '  +6.89%  '

$ws.Range("E42").Value = '  -0.26%  '

$ws.Range("E43").Value = '  -1.96%  '

$ws.Range("E44").Value = '  -0.31%  '

$ws.Range("E45").Value = '  +0.18%  '

$ws.Range("E46").Value = '  -0.80%  '

Set-TextValue $ws.Range("D47") '1.782.57'
$ws.Range("E47").Value = '  +0.06%  '

$ws.Range("E48").Value = '  +2.51%  '

Set-TextValue $ws.Range("D49") '87.09'
$ws.Range("E49").Value = '  +1.32%  '

$ws.Range("E50").Value = '  -0.35%  '

$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range("D51") '0.0993'
$ws.Range("E51").Value = '  -0.01%  '
